$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (D new value or $null, E new value or $null)
# D values that parse as plain numbers must be forced to stay text
# (matching the original price-column formatting, e.g. '24.60' not 24.6),
# so NumberFormat is temporarily set to Text, then cleared again so the
# cell keeps the default/general style like the rest of the sheet.
$updates = @(
    @{ Row = 2; D = "30.616.48"; E = "  +0.61%  " }
    @{ Row = 3; D = "2.116.13"; E = "  +1.17%  " }
    @{ Row = 4; D = "1.011"; E = "  +0.90%  " }
    @{ Row = 5; D = "336.87"; E = "  +2.02%  " }
    @{ Row = 6; D = $null; E = "  +0.84%  " }
    @{ Row = 7; D = "0.5241"; E = "  +0.48%  " }
    @{ Row = 8; D = "0.4556"; E = "  +3.99%  " }
    @{ Row = 9; D = "54.48"; E = "  +1.48%  " }
    @{ Row = 10; D = "0.09106"; E = "  +2.44%  " }
    @{ Row = 11; D = "1.172"; E = "  +1.51%  " }
    @{ Row = 12; D = "24.60"; E = "  +1.30%  " }
    @{ Row = 13; D = "2.113.84"; E = "  +1.33%  " }
    @{ Row = 14; D = "6.862"; E = "  +2.41%  " }
    @{ Row = 15; D = "8.089"; E = "  +5.22%  " }
    @{ Row = 16; D = "0.00001174"; E = "  +4.73%  " }
    @{ Row = 17; D = "97.01"; E = "  +1.12%  " }
    @{ Row = 18; D = "1.011"; E = "  +0.94%  " }
    @{ Row = 19; D = $null; E = "  +1.43%  " }
    @{ Row = 20; D = "19.41"; E = "  +1.17%  " }
    @{ Row = 22; D = "6.301"; E = "  +0.70%  " }
    @{ Row = 23; D = "30.646.64"; E = "  +0.60%  " }
    @{ Row = 24; D = "12.86"; E = "  +4.89%  " }
    @{ Row = 25; D = "2.360"; E = "  +1.34%  " }
    @{ Row = 26; D = "2.366.09"; E = "  +1.47%  " }
    @{ Row = 27; D = "22.34"; E = "  +0.41%  " }
    @{ Row = 28; D = "163.81"; E = "  +0.65%  " }
    @{ Row = 29; D = "2.539"; E = "  -0.79%  " }
    @{ Row = 30; D = "134.72"; E = "  +2.48%  " }
    @{ Row = 31; D = "1.207"; E = "  +1.95%  " }
    @{ Row = 33; D = "1.653"; E = "  -0.96%  " }
    @{ Row = 34; D = "6.374"; E = "  +3.32%  " }
    @{ Row = 35; D = "3.955"; E = "  +1.53%  " }
    @{ Row = 36; D = "10.57"; E = "  +5.24%  " }
    @{ Row = 37; D = "5.943"; E = "  +8.47%  " }
    @{ Row = 38; D = "0.02623"; E = "  +2.27%  " }
    @{ Row = 39; D = "0.06836"; E = "  +0.14%  " }
    @{ Row = 40; D = $null; E = "  +3.18%  " }
    @{ Row = 41; D = $null; E = "  -0.25%  " }
    @{ Row = 42; D = "0.6872"; E = "  -0.11%  " }
    @{ Row = 43; D = "1.260"; E = "  +0.57%  " }
    @{ Row = 44; D = "14.87"; E = "  +6.25%  " }
    @{ Row = 45; D = "0.6458"; E = "  +1.79%  " }
    @{ Row = 46; D = "2.316"; E = "  +5.31%  " }
    @{ Row = 47; D = $null; E = "  +22.47%  " }
    @{ Row = 48; D = "3.687"; E = "  +1.81%  " }
    @{ Row = 49; D = $null; E = "  +0.85%  " }
    @{ Row = 50; D = "83.59"; E = "  +2.22%  " }
    @{ Row = 51; D = "0.3355"; E = "  +12.61%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)   # column D = Price
        if ($u.D -match '^-?\d+(\.\d+)?$') {
            # Value reads as a plain number (e.g. "24.60") - Excel would silently
            # coerce it (and drop the trailing zero) unless the cell is forced to
            # Text first. Clear the format again afterwards so the cell ends up
            # with the same (default) style as before, only the content differs.
            $cell.NumberFormat = "@"
            $cell.Value = $u.D
            $cell.ClearFormats()
        } else {
            $cell.Value = $u.D
        }
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E   # column E = Volume(1h)
    }
}
